$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" (column G) values were re-ordered (same set of names/emails,
# different order within the comma-separated list) by the daily attendance
# processing job. Apply the same canonical re-ordering here.
$map = @{
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $cur = $cell.Value2
    if ($cur -ne $null -and $map.ContainsKey($cur)) {
        $cell.Value = $map[$cur]
    }
}

Write-Host "DONE"
